# Conserto do erro com o rotulo da coluna 2050 nas tabelas
# e retirada das linhas com total das tabelas
$wb = $excel.ActiveWorkbook

# Sheet 1: Potencia Acumulada - SIN (MW)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("E1").Value = "'2050"
$ws1.Rows.Item(13).Delete()

# Sheet 2: Geracao Periodo Medio (MWMed)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("E1").Value = "'2050"
$ws2.Rows.Item(13).Delete()

# Sheet 3: Atendimento a Ponta(MW)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("E1").Value = "'2050"
$ws3.Rows.Item(13).Delete()

# Sheet 4: Potencia Incremental - SIN(MW) -- header uses a range label (2041-2050)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("E1").Value = "'2041-2050"
$ws4.Rows.Item(13).Delete()

# Sheet 5: Emissoes Totais (MtCO2eq) -- only the label fix, no Total row present
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("E1").Value = "'2050"

# Sheet 6: Custo Total (bilhoes de R$) -- only the Total row removal
$ws6 = $wb.Worksheets.Item(6)
$ws6.Rows.Item(4).Delete()
